$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo/casing in the "admin staff" role label used for group 10 (row 2).
$ws.Range("B2").Value = "Department of System and Computer Engineering - Admin Staff "

# Fill in the (previously empty) "user access role" column for each project group.
$ws.Range("D2").Value = "Ish"
$ws.Range("D7").Value = "-"

# Mark the assessmentStatus flag as complete for Russ, Daren (row 8) and Singh, Ishdeep (row 10).
$ws.Range("H8").Value = 1
$ws.Range("H10").Value = 1

# Update the saved selection to match the latest interaction (the merged
# "user access role" block for the first group).
$ws.Range("D7:D10").Select()
